# Regenerate the s_vals data (TB, d2S, K, IP columns and the derived sum)
# to reflect filtering out "save games" from the underlying computation.
# Column F (Win) is untouched; column G (sum) = B + C + D + E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row: B (TB), C (d2S), D (K), E (IP)
$data = @{
    2  = @(3.230985683306322,  1.667794583268128,  3.900430680208489,  0.496779210170732)
    3  = @(3.230985683306322,  1.667794583268128,  3.900430680208489,  0.496779210170732)
    4  = @(0.127881588408715,  0.3127903958511391, 26.21740644021617,  0.496779210170732)
    5  = @(1.459612070389937,  1.667794583268128,  0.8054896365839992, 0.496779210170732)
    6  = @(1.459612070389937,  1.667794583268128,  3.900430680208489,  8.660232485948974)
    7  = @(0.6753301551942219, 1.667794583268128,  3.900430680208489,  0.496779210170732)
    8  = @(3.230985683306322,  1.667794583268128,  0.1575252929769615, 0.496779210170732)
    9  = @(0.04763786555579896,0.3127903958511391, 0.1575252929769615, 0.496779210170732)
    10 = @(0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 0.496779210170732)
    11 = @(0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 0.496779210170732)
    12 = @(1.459612070389937,  1.667794583268128,  3.900430680208489,  0.496779210170732)
    13 = @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 7).Value = $b + $c + $d + $e
}
